$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Teams")

# Insert a new column for "Channel1Description" right after "Channel1Name" (old column E)
$ws.Columns("E").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)

# Insert a new column for "Channel2Description" right after "Channel2Name" (old column F, now G after first insert)
$ws.Columns("H").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)

# Headers and data, entered in the same order as the authored workbook so
# that the shared-strings table is built up identically.
$ws.Range("E1").Value = "Channel1Description"
$ws.Range("E2").Value = "Channel for all Marketing Campaigns"
$ws.Range("H1").Value = "Channel2Description"
$ws.Range("H2").Value = "Channel for Brainstorming"
$ws.Range("H4").Value = "Channel for all agreements"
$ws.Range("H5").Value = "A private channel for confidential information"
$ws.Range("H3").Value = "A private channel for projects"
$ws.Range("E3").Value = "Channel for all Incidents"
$ws.Range("E4").Value = "Channel for all contracts"
$ws.Range("E5").Value = "Channel for all onboarding"

# Update selection to match the target state
$ws.Range("E6").Select()
